$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the daily figures for the rows that were revised in this upload ---
# Row 280 (2020-12-01): confirmed cases revised 125 -> 126
$ws.Range("C280").Value = 126

# Row 281 (2020-12-02): confirmed cases revised 94 -> 117, deaths detail 2/0 -> 3/1
$ws.Range("C281").Value = 117
$ws.Range("L281").Value = 3
$ws.Range("M281").Value = 1

# Row 282 (2020-12-03): confirmed cases revised 19 -> 85
$ws.Range("C282").Value = 85

# Row 283 (2020-12-04): row now filled in with real data instead of being blank
$ws.Range("C283").Value = 22
$ws.Range("E283").Value = 20
$ws.Range("F283").Value = 20
$ws.Range("G283").Value = 119
$ws.Range("L283").Value = 0
$ws.Range("M283").Value = 0

# --- Restore the view state (scroll position / selection) recorded in the sheet ---
$excel.ActiveWindow.ScrollRow = 256
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("L3:M283").Select()
